# Insert a new price-history column (AH) before the existing "nom" (AH) /
# "url_produit" (AI) columns, shifting them one column to the right
# (to AI / AJ respectively), then populate the new column for every
# product row that still has a price (rows 2-80 carry the same price
# as the previous column AG; rows 81-205 have no price and stay blank,
# matching the existing "no data" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 34 = AH (A=1 ... Z=26, AA=27 ... AH=34)
$ws.Columns.Item(34).Insert()

# New header timestamp for the freshly inserted column.
$ws.Range("AH1").Value = "2026-01-29 05:44:42"

# Carry forward the last known price (from column AG, now column AG
# still since it is to the left of the insertion point) into the new
# AH column for each product row that has a price.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 33).Value2
    $ws.Cells.Item($r, 34).Value = $price
}
